# Regenerate save_data to use K (strikeouts) instead of Strike# in column G.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K values for rows 2-38 (column G), replacing the old Strike# values.
$kValues = @{
    2  = 3
    3  = 0
    4  = 1
    5  = 3
    6  = 3
    7  = 5
    8  = 2
    9  = 2
    10 = 0
    11 = 7
    12 = 4
    13 = 4
    14 = 6
    15 = 5
    16 = 3
    17 = 8
    18 = 4
    19 = 4
    20 = 4
    21 = 2
    22 = 4
    23 = 5
    24 = 5
    25 = 8
    26 = 5
    27 = 4
    28 = 4
    29 = 5
    30 = 6
    31 = 6
    32 = 2
    33 = 4
    34 = 4
    35 = 5
    36 = 1
    37 = 4
    38 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
